$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 24000.2
$ws.Range("I6").Value = 5667
$ws.Range("J6").Value = 51500
$ws.Range("K6").Value = 17001
$ws.Range("L6").Value = 154500
$ws.Range("M6").Value = -16889
$ws.Range("N6").Value = -154724
$ws.Range("H33").Value = 258.22974
$ws.Range("I33").Value = 197.30986
$ws.Range("K33").Value = 197.30986
$ws.Range("M33").Value = 31.69014000000001
$ws.Range("H64").Value = 3200
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3200
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 3200
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3200
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4916
$ws.Range("H70").Value = 1571.4286
$ws.Range("I70").Value = 1666.6666
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 4999.9998
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -4729.9998
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 1571.4286
$ws.Range("I73").Value = 1666.6666
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 4999.9998
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -4063.9998
$ws.Range("N73").Value = -6372
$ws.Range("H76").Value = 3131.1035
$ws.Range("I76").Value = 3002.4
$ws.Range("J76").Value = 3198.842
$ws.Range("K76").Value = 3002.4
$ws.Range("L76").Value = 3198.842
$ws.Range("M76").Value = -2687.4
$ws.Range("N76").Value = -3828.842
$ws.Range("H79").Value = 3131.1035
$ws.Range("I79").Value = 3002.4
$ws.Range("J79").Value = 3198.842
$ws.Range("K79").Value = 3002.4
$ws.Range("L79").Value = 3198.842
$ws.Range("M79").Value = -1910.4
$ws.Range("N79").Value = -5382.842000000001
$ws.Range("H97").Value = 551000
$ws.Range("J97").Value = 551000
$ws.Range("L97").Value = 1653000
$ws.Range("N97").Value = -1653992
$ws.Range("H99").Value = 3352.8
$ws.Range("I99").Value = 3588
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 10764
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -9266
$ws.Range("N99").Value = -11996
$ws.Range("H100").Value = 2374.875
$ws.Range("I100").Value = 1399.8
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 1399.8
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -858.8
$ws.Range("N100").Value = -5082
$ws.Range("H112").Value = 1472421.2
$ws.Range("I112").Value = 866.6667
$ws.Range("J112").Value = 1787754.4
$ws.Range("K112").Value = 2600.0001
$ws.Range("L112").Value = 5363263.199999999
$ws.Range("M112").Value = -1492.0001
$ws.Range("N112").Value = -5365479.199999999
$ws.Range("H117").Value = 48723.75
$ws.Range("J117").Value = 48723.75
$ws.Range("L117").Value = 48723.75
$ws.Range("N117").Value = -57901.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H117").Value = 49561.75
$ws.Range("J117").Value = 49561.75
$ws.Range("L117").Value = 49561.75
$ws.Range("N117").Value = -58739.75
$ws.Range("H118").Value = 49803
$ws.Range("J118").Value = 49803
$ws.Range("L118").Value = 49803
$ws.Range("N118").Value = -53117
$ws.Range("H135").Value = 50376.2
$ws.Range("J135").Value = 50376.2
$ws.Range("L135").Value = 50376.2
$ws.Range("N135").Value = -60516.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 27875
$ws.Range("J33").Value = 36333.332
$ws.Range("L33").Value = 36333.332
$ws.Range("N33").Value = -37005.332
$ws.Range("H37").Value = 2045.2
$ws.Range("I37").Value = 1113
$ws.Range("K37").Value = 1113
$ws.Range("M37").Value = -976
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 43468
$ws.Range("J112").Value = 43468
$ws.Range("L112").Value = 43468
$ws.Range("N112").Value = -46422
$ws.Range("H132").Value = 34888.215
$ws.Range("I132").Value = 1152.8235
$ws.Range("J132").Value = 178263.62
$ws.Range("K132").Value = 3458.4705
$ws.Range("L132").Value = 534790.86
$ws.Range("M132").Value = -928.4704999999999
$ws.Range("N132").Value = -539850.86
$ws.Range("H134").Value = 242562.23
$ws.Range("I134").Value = 1074.9608
$ws.Range("J134").Value = 2001969.6
$ws.Range("K134").Value = 3224.8824
$ws.Range("L134").Value = 6005908.800000001
$ws.Range("M134").Value = -689.8824000000004
$ws.Range("N134").Value = -6010978.800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 10000280
$ws.Range("I2").Value = 41.11111
$ws.Range("J2").Value = 15625414
$ws.Range("K2").Value = 246.66666
$ws.Range("L2").Value = 93752484
$ws.Range("M2").Value = -133.66666
$ws.Range("N2").Value = -93752710
$ws.Range("H3").Value = 6712.8
$ws.Range("I3").Value = 2978.3333
$ws.Range("J3").Value = 8313.286
$ws.Range("K3").Value = 8934.999899999999
$ws.Range("L3").Value = 24939.858
$ws.Range("M3").Value = -8822.999899999999
$ws.Range("N3").Value = -25163.858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 505
$ws.Range("I19").Value = 505
$ws.Range("K19").Value = 505
$ws.Range("M19").Value = -217
$ws.Range("H110").Value = 49233.332
$ws.Range("J110").Value = 49233.332
$ws.Range("L110").Value = 49233.332
$ws.Range("N110").Value = -57413.332
$ws.Range("H136").Value = 63442
$ws.Range("J136").Value = 63442
$ws.Range("L136").Value = 190326
$ws.Range("N136").Value = -195426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5220.4736
$ws.Range("I46").Value = 1918.9
$ws.Range("J46").Value = 8888.888999999999
$ws.Range("K46").Value = 1918.9
$ws.Range("L46").Value = 8888.888999999999
$ws.Range("M46").Value = -1730.9
$ws.Range("N46").Value = -9264.888999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 16527.309
$ws.Range("J70").Value = 16527.309
$ws.Range("L70").Value = 16527.309
$ws.Range("N70").Value = -17157.309
$ws.Range("H73").Value = 16527.309
$ws.Range("J73").Value = 16527.309
$ws.Range("L73").Value = 16527.309
$ws.Range("N73").Value = -18711.309
$ws.Range("H136").Value = 233371.11
$ws.Range("I136").Value = 286317.5
$ws.Range("J136").Value = 1730.75
$ws.Range("K136").Value = 858952.5
$ws.Range("L136").Value = 5192.25
$ws.Range("M136").Value = -856402.5
$ws.Range("N136").Value = -10292.25
